$d = $word.ActiveDocument

$replacements = @(
    @{old = "613÷7=87, 4"; new = "281÷3=93, 2"},
    @{old = "720÷9=80, 0"; new = "295÷7=42, 1"},
    @{old = "730÷2=365, 0"; new = "568÷7=81, 1"},
    @{old = "383÷7=54, 5"; new = "866÷9=96, 2"},
    @{old = "331÷3=110, 1"; new = "377÷2=188, 1"},
    @{old = "789÷8=98, 5"; new = "477÷2=238, 1"},
    @{old = "279÷4=69, 3"; new = "922÷2=461, 0"},
    @{old = "383÷2=191, 1"; new = "824÷3=274, 2"},
    @{old = "369÷3=123, 0"; new = "709÷2=354, 1"},
    @{old = "760÷8=95, 0"; new = "144÷8=18, 0"},
    @{old = "231÷2=115, 1"; new = "642÷8=80, 2"},
    @{old = "397÷6=66, 1"; new = "668÷6=111, 2"},
    @{old = "549÷7=78, 3"; new = "369÷2=184, 1"},
    @{old = "857÷4=214, 1"; new = "836÷7=119, 3"},
    @{old = "808÷6=134, 4"; new = "930÷9=103, 3"},
    @{old = "512÷7=73, 1"; new = "949÷2=474, 1"},
    @{old = "107÷8=13, 3"; new = "889÷6=148, 1"},
    @{old = "490÷6=81, 4"; new = "430÷3=143, 1"},
    @{old = "320÷8=40, 0"; new = "616÷3=205, 1"},
    @{old = "372÷6=62, 0"; new = "198÷6=33, 0"},
    @{old = "976÷5=195, 1"; new = "752÷8=94, 0"},
    @{old = "708÷4=177, 0"; new = "305÷4=76, 1"},
    @{old = "494÷3=164, 2"; new = "719÷3=239, 2"},
    @{old = "722÷7=103, 1"; new = "514÷4=128, 2"},
    @{old = "442÷9=49, 1"; new = "713÷8=89, 1"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
